$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: fill in the rest of the week's data (C, D updated, E, J, K added; G/H/I updated)
$ws.Range("C46").Value = 94.56999999999999
$ws.Range("D46").Value = 248.8
$ws.Range("E46").Value = 154.23
$ws.Range("G46").Value = 1.67
$ws.Range("H46").Value = 206.83
$ws.Range("I46").Value = 255.29
$ws.Range("J46").Value = 164.72
$ws.Range("K46").Value = 159.18

# Row 47: the prediction value is revised
$ws.Range("D47").Value = 189.49

# Row 48: new week of prediction data (previously row 46's partial data)
# The leading apostrophe forces text so Excel doesn't auto-convert the
# date-shaped string into a date serial; Style reset clears the
# quote-prefix formatting that Excel applies to such cells.
$ws.Range("A48").Value = "'2021-01-09"
$ws.Range("A48").Style = "Normal"
$ws.Range("B48").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D48").Value = 175.16
$ws.Range("F48").Value = "KNN"
$ws.Range("G48").Value = 1.02
$ws.Range("H48").Value = 128.61
$ws.Range("I48").Value = 121.63

# Row 49: new week of prediction data (previously row 47's data)
$ws.Range("A49").Value = "'2021-01-09"
$ws.Range("A49").Style = "Normal"
$ws.Range("B49").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D49").Value = 128.33
$ws.Range("F49").Value = "KNN"
